$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.355.31'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').Value = '2.603.75'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.34'
$ws.Range('E5').Value = '  +6.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.95'
$ws.Range('E6').Value = '  +1.16%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.599'
$ws.Range('E8').Value = '  +0.58%  '
$ws.Range('D9').Value = '2.614.19'
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.51'
$ws.Range('E10').Value = '  -3.34%  '
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('E13').Value = '  +4.20%  '
$ws.Range('D14').Value = '3.068.23'
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.68'
$ws.Range('E15').Value = '  +6.14%  '
$ws.Range('D16').Value = '60.366.68'
$ws.Range('E16').Value = '  +1.20%  '
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').Value = '2.612.17'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('E19').Value = '  +9.28%  '
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '347.63'
$ws.Range('E21').Value = '  +2.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.91'
$ws.Range('E22').Value = '  +5.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E24').Value = '  +8.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.12'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.160'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.09'
$ws.Range('E28').Value = '  +8.01%  '
$ws.Range('E29').Value = '  +2.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.88'
$ws.Range('E30').Value = '  +11.05%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.38'
$ws.Range('E31').Value = '  +3.18%  '
$ws.Range('B32').Value = 'USDe'
$ws.Range('C32').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.998'
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '163.38'
$ws.Range('E33').Value = '  +3.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.49'
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.28'
$ws.Range('E35').Value = '  +3.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.980'
$ws.Range('E36').Value = '  +7.86%  '
$ws.Range('E37').Value = '  +5.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.63'
$ws.Range('E38').Value = '  +9.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.07'
$ws.Range('E39').Value = '  +1.11%  '
$ws.Range('E40').Value = '  +6.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '310.83'
$ws.Range('E41').Value = '  +6.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.838'
$ws.Range('E42').Value = '  -0.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '135.89'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0995'
$ws.Range('E44').Value = '  +2.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.80'
$ws.Range('E46').Value = '  +3.49%  '
$ws.Range('E47').Value = '  +5.36%  '
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('E49').Value = '  +2.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.14'
$ws.Range('E50').Value = '  +7.47%  '
$ws.Range('E51').Value = '  +2.58%  '
